$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07633962631225585
$ws.Range("C2").Value = 0.02074279764011377
$ws.Range("D2").Value = 0.2543339729309082
$ws.Range("E2").Value = 0.01035063928488965

$ws.Range("B3").Value = 0.06673617362976074
$ws.Range("C3").Value = 0.002448247433804212
$ws.Range("D3").Value = 0.294121265411377
$ws.Range("E3").Value = 0.0691875972842113

$ws.Range("B4").Value = 0.06385512351989746
$ws.Range("C4").Value = 0.001170061260147135
$ws.Range("D4").Value = 0.2036287784576416
$ws.Range("E4").Value = 0.008945215755263723

$ws.Range("B5").Value = 0.07987580299377442
$ws.Range("C5").Value = 0.0184100325933261
$ws.Range("D5").Value = 0.234155797958374
$ws.Range("E5").Value = 0.025029728461905

$ws.Range("B6").Value = 0.08663024902343749
$ws.Range("C6").Value = 0.02480285007467933
$ws.Range("D6").Value = 0.4677600860595703
$ws.Range("E6").Value = 0.1185483050510801

$ws.Range("B7").Value = 0.09374814033508301
$ws.Range("C7").Value = 0.01646580248803083
$ws.Range("D7").Value = 0.4695443153381348
$ws.Range("E7").Value = 0.0853797996697889

$ws.Range("B8").Value = 0.07477412223815919
$ws.Range("C8").Value = 0.01018735420539713
$ws.Range("D8").Value = 0.3012666702270508
$ws.Range("E8").Value = 0.05821066291779752

$ws.Range("B9").Value = 0.07215652465820313
$ws.Range("C9").Value = 0.008530893259487735
$ws.Range("D9").Value = 0.2534438610076905
$ws.Range("E9").Value = 0.01644369689215851

$wb.Save()
